# Updated cryptos list on Mon Apr 15 18:48:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as TEXT (matching the workbook's inline-string
# cells) without letting Excel's type-inference turn numeric-looking
# strings (e.g. "554.21", "1.00", "0.0410") into real numbers. The
# NumberFormat is restored to General afterwards so no stray formatting
# is left behind.
function Set-TextCell {
    param($cellRef, $value)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.NumberFormat = "General"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "63.751.43"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.087.63"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextCell "D5" "554.21"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6 - Solana
Set-TextCell "D6" "137.31"
$ws.Range("E6").Value = "  -3.76%  "

# Row 7 - USDC
Set-TextCell "D7" "0.999"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - LidoStakedEther
Set-TextCell "D8" "3.085.52"
$ws.Range("E8").Value = "  +0.23%  "

# Row 9 - XRP
Set-TextCell "D9" "0.493"
$ws.Range("E9").Value = "  +0.12%  "

# Row 10 - Toncoin
Set-TextCell "D10" "6.65"
$ws.Range("E10").Value = "  +1.07%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.156"
$ws.Range("E11").Value = "  +1.61%  "

# Row 12 - Cardano
Set-TextCell "D12" "0.449"
$ws.Range("E12").Value = "  -1.98%  "

# Row 13 - Avalanche
Set-TextCell "D13" "34.96"
$ws.Range("E13").Value = "  -3.51%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -1.67%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell "D15" "3.582.72"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "63.749.05"
$ws.Range("E16").Value = "  -1.16%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.24%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "3.084.70"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19 - BitcoinCash
Set-TextCell "D19" "504.53"
$ws.Range("E19").Value = "  +2.56%  "

# Row 20 - Polkadot
Set-TextCell "D20" "6.59"
$ws.Range("E20").Value = "  -1.22%  "

# Row 21 - Chainlink
Set-TextCell "D21" "13.52"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.02%  "

# Row 23 - Uniswap
Set-TextCell "D23" "7.17"
$ws.Range("E23").Value = "  -1.76%  "

# Rows 24/25 swap: Litecoin <-> InternetComputer(DFINITY)
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D24" "12.29"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "77.04"
$ws.Range("E25").Value = "  -2.12%  "

# Row 26 - Dai
Set-TextCell "D26" "0.999"

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  +0.85%  "

# Row 28 - RenderToken
Set-TextCell "D28" "8.23"
$ws.Range("E28").Value = "  +1.39%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  -3.56%  "

# Row 30 - FirstDigitalUSD
Set-TextCell "D30" "1.00"
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - EthereumClassic
Set-TextCell "D31" "26.01"
$ws.Range("E31").Value = "  -0.13%  "

# Row 32 - Stacks
Set-TextCell "D32" "2.52"
$ws.Range("E32").Value = "  -5.38%  "

# Row 33 - Mantle
Set-TextCell "D33" "1.11"
$ws.Range("E33").Value = "  -3.02%  "

# Row 34 - Bittensor
Set-TextCell "D34" "530.45"
$ws.Range("E34").Value = "  -11.42%  "

# Row 35 - OKB
Set-TextCell "D35" "57.93"
$ws.Range("E35").Value = "  +10.74%  "

# Row 36 - Filecoin
Set-TextCell "D36" "5.84"
$ws.Range("E36").Value = "  -3.07%  "

# Row 37 - NEARProtocol
Set-TextCell "D37" "5.13"
$ws.Range("E37").Value = "  -6.25%  "

# Row 38 - VeChain
Set-TextCell "D38" "0.0410"
$ws.Range("E38").Value = "  +1.14%  "

# Row 39 - Hedera
Set-TextCell "D39" "0.0791"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40 - Maker
Set-TextCell "D40" "3.045.33"
$ws.Range("E40").Value = "  +2.16%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -1.70%  "

# Row 42 - Cosmos
Set-TextCell "D42" "8.05"
$ws.Range("E42").Value = "  -2.77%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -11.31%  "

# Row 45 - TheGraph
Set-TextCell "D45" "0.249"
$ws.Range("E45").Value = "  +0.40%  "

# Row 46 - Fetch.AI
$ws.Range("E46").Value = "  -3.77%  "

# Row 47 - Monero
Set-TextCell "D47" "122.19"
$ws.Range("E47").Value = "  +1.51%  "

# Rows 48/49 swap: Stellar <-> InjectiveProtocol
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D48" "23.88"
$ws.Range("E48").Value = "  -5.25%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D49" "0.106"
$ws.Range("E49").Value = "  -1.76%  "

# Row 50 - PEPE
Set-TextCell "D50" "0.0₃0491"
$ws.Range("E50").Value = "  -8.98%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  -4.14%  "
